$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.639.16"
$ws.Range("E2").Value = "  +5.84%  "
$ws.Range("D3").Value = "2.744.12"
$ws.Range("E3").Value = "  +4.45%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "332.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.533"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +6.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0832"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.95%  "
$ws.Range("D15").Value = "3.171.33"
$ws.Range("E15").Value = "  +4.50%  "
$ws.Range("D16").Value = "2.735.83"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").Value = "51.519.92"
$ws.Range("E18").Value = "  +5.66%  "
$ws.Range("B19").Value = "ImmutableX"
$ws.Range("C19").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.56%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.02%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0820"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.39%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0346"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.47%  "
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("E44").Value = "  +4.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.50%  "
$ws.Range("D46").Value = "2.099.46"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.52%  "
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.51%  "
